$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 0
